$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.540008666666668
$ws.Range("H2").Value = 19.620026
$ws.Range("I2").Value = 0.2365207520404831
$ws.Range("J2").Value = 0.2365207520404831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 112.2692549765933
$ws.Range("R2").Value = 1010.42329478934
$ws.Range("S2").Value = 0.01325333147312037
$ws.Range("T2").Value = 0.01325333147312037

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.540008666666668
$ws.Range("H3").Value = 19.620026
$ws.Range("I3").Value = 0.2365207520404831
$ws.Range("J3").Value = 0.2365207520404831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 1677.147979397313
$ws.Range("R3").Value = 15094.33181457582
$ws.Range("S3").Value = 0.1979865111339773
$ws.Range("T3").Value = 0.1979865111339773

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.540008666666668
$ws.Range("H4").Value = 19.620026
$ws.Range("I4").Value = 0.2365207520404831
$ws.Range("J4").Value = 0.2365207520404831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 214.1551256733685
$ws.Range("R4").Value = 1927.396131060316
$ws.Range("S4").Value = 0.02528090943338534
$ws.Range("T4").Value = 0.02528090943338534

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.67485766666667
$ws.Range("H5").Value = 44.024573
$ws.Range("I5").Value = 0.5307192311682535
$ws.Range("J5").Value = 0.5307192311682536
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 251.9163843805633
$ws.Range("R5").Value = 2267.24745942507
$ws.Range("S5").Value = 0.02973860783525899
$ws.Range("T5").Value = 0.02973860783525899

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.67485766666667
$ws.Range("H6").Value = 44.024573
$ws.Range("I6").Value = 0.5307192311682535
$ws.Range("J6").Value = 0.5307192311682536
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 3763.283680193874
$ws.Range("R6").Value = 33869.55312174487
$ws.Range("S6").Value = 0.4442538257815304
$ws.Range("T6").Value = 0.4442538257815305

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.67485766666667
$ws.Range("H7").Value = 44.024573
$ws.Range("I7").Value = 0.5307192311682535
$ws.Range("J7").Value = 0.5307192311682536
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 480.5339179230131
$ws.Range("R7").Value = 4324.805261307119
$ws.Range("S7").Value = 0.05672679755146406
$ws.Range("T7").Value = 0.05672679755146406

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 6.436020999999999
$ws.Range("H8").Value = 19.308063
$ws.Range("I8").Value = 0.2327600167912634
$ws.Range("J8").Value = 0.2327600167912634
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 110.48414757713
$ws.Range("R8").Value = 994.3573281941698
$ws.Range("S8").Value = 0.01304260040444854
$ws.Range("T8").Value = 0.01304260040444854

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 6.436020999999999
$ws.Range("H9").Value = 19.308063
$ws.Range("I9").Value = 0.2327600167912634
$ws.Range("J9").Value = 0.2327600167912634
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 1650.480934455745
$ws.Range("R9").Value = 14854.32841010171
$ws.Range("S9").Value = 0.1948384793233727
$ws.Range("T9").Value = 0.1948384793233727

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 6.436020999999999
$ws.Range("H10").Value = 19.308063
$ws.Range("I10").Value = 0.2327600167912634
$ws.Range("J10").Value = 0.2327600167912634
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 210.7500091118286
$ws.Range("R10").Value = 1896.750082006458
$ws.Range("S10").Value = 0.02487893706344213
$ws.Range("T10").Value = 0.02487893706344213
